$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 2.42
$ws.Range("K2").Value = 4.7
$ws.Range("O2").Value = 1.21
$ws.Range("Q2").Value = 1.65
$ws.Range("W2").Value = 1.53
# Row 3 updates
$ws.Range("I3").Value = 5.9
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 5.1
$ws.Range("P3").Value = 2.8
$ws.Range("Q3").Value = 1.47
$ws.Range("R3").Value = 1.74
$ws.Range("S3").Value = 2.16
$ws.Range("U3").Value = 2.18
$ws.Range("V3").Value = 1.21
$ws.Range("Y3").Value = 980
$ws.Range("Z3").Value = 55
$ws.Range("AD3").Value = 23
$ws.Range("AF3").Value = 13
$ws.Range("AH3").Value = 17.5
$ws.Range("AJ3").Value = 17.5
$ws.Range("AO3").Value = 44
# Row 4 updates
$ws.Range("F4").Value = 2.76
$ws.Range("H4").Value = 2.46
$ws.Range("I4").Value = 2.72
$ws.Range("K4").Value = 4
$ws.Range("N4").Value = 4.3
$ws.Range("O4").Value = 1.24
$ws.Range("Q4").Value = 1.74
$ws.Range("R4").Value = 1.46
$ws.Range("S4").Value = 2.6
$ws.Range("T4").Value = 1.61
$ws.Range("U4").Value = 2.32
$ws.Range("V4").Value = 1.58
$ws.Range("W4").Value = 1.49
$ws.Range("X4").Value = 20
$ws.Range("Y4").Value = 16
$ws.Range("Z4").Value = 23
$ws.Range("AA4").Value = 44
$ws.Range("AB4").Value = 16.5
$ws.Range("AC4").Value = 10.5
$ws.Range("AD4").Value = 14.5
$ws.Range("AE4").Value = 32
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 16
$ws.Range("AH4").Value = 19
$ws.Range("AI4").Value = 42
$ws.Range("AJ4").Value = 55
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 980
$ws.Range("AM4").Value = 85
$ws.Range("AN4").Value = 27
$ws.Range("AO4").Value = 22
# Row 5 updates
$ws.Range("F5").Value = 3.5
$ws.Range("G5").Value = 3.75
$ws.Range("H5").Value = 2.08
$ws.Range("I5").Value = 2.14
$ws.Range("J5").Value = 3.85
$ws.Range("K5").Value = 4.2
$ws.Range("N5").Value = 5.7
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 2.58
$ws.Range("Q5").Value = 1.54
$ws.Range("R5").Value = 1.64
$ws.Range("S5").Value = 2.3
$ws.Range("T5").Value = 1.53
$ws.Range("U5").Value = 2.64
$ws.Range("V5").Value = 1.88
$ws.Range("W5").Value = 1.36
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 18.5
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 980
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 12
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 16
$ws.Range("AH5").Value = 18
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 60
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 55
$ws.Range("AO5").Value = 9.6
# Row 6 updates
$ws.Range("F6").Value = 7.2
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 1.5
$ws.Range("J6").Value = 4.8
$ws.Range("L6").Value = 1.3
$ws.Range("N6").Value = 5.1
$ws.Range("R6").Value = 1.58
$ws.Range("V6").Value = 2.92
$ws.Range("Z6").Value = 12.5
$ws.Range("AC6").Value = 13.5
$ws.Range("AD6").Value = 13
$ws.Range("AJ6").Value = 210

$wb.Save()
